$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a Text number format on cells whose new values are plain numbers,
# so Excel stores them as text (matching the source inlineStr cells) instead
# of auto-converting them to numeric cells. The format is cleared back to the
# default "Normal" style right after the value is written so the cells keep no
# explicit style attribute, just like the rest of the sheet.
$numericTextCells = @("D4", "D5", "D7", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D20", "D22", "D23", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.254.55"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.861.80"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "236.30"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.4716"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("E8").Value = "  +2.29%  "
$ws.Range("D9").Value = "0.06550"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "21.90"
$ws.Range("E10").Value = "  +2.87%  "
$ws.Range("D11").Value = "0.07949"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "97.70"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "1.864.03"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "5.142"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "0.6799"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "265.14"
$ws.Range("E16").Value = "  -5.46%  "
$ws.Range("D17").Value = "30.248.78"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("E18").Value = "  +7.98%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "0.000007494"
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("D21").Value = "2.109.17"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "5.278"
$ws.Range("E23").Value = "  -4.11%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "167.83"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("D26").Value = "9.189"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").Value = "1.947"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "1.395"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").Value = "0.09959"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("D31").Value = "4.341"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("D32").Value = "1.469"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").Value = "4.015"
$ws.Range("E33").Value = "  -1.97%  "
$ws.Range("D34").Value = "0.04703"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").Value = "0.7005"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").Value = "2.625"
$ws.Range("E39").Value = "  +3.49%  "
$ws.Range("D40").Value = "6.309"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").Value = "73.83"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "1.934"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8409"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").Value = "0.4160"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "103.32"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("D47").Value = "7.144"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").Value = "944.14"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").Value = "9.169"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "34.15"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "0.05663"
$ws.Range("E51").Value = "  +0.65%  "

foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}
